$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the two runs of the last paragraph ("Eltern nicht fertig" /
#    " Do File fertigstellen") into a single run and drop the
#    "_GoBack" bookmark that used to sit between them. A Find/Replace
#    over the full text re-writes the range as one run and removes the
#    bookmark anchor that fell inside the replaced span.
# ------------------------------------------------------------------
$old = "„Eltern nicht fertig“ Do File fertigstellen"
$new = "„Eltern nicht fertig“ Do File fertigstellen"
$found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Safety net: if, for some reason, the bookmark survived the replace
# above (e.g. it had moved to the very end), make sure it is gone
# before we re-add it after the new paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Insert the new to-do item as its own list paragraph right after
#    the paragraph we just normalised, carrying over the bookmark.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
# Land just before the paragraph mark (Range.End - 1) rather than at
# the absolute end of the story: inserting exactly at Content.End
# overwrites the trailing paragraph instead of adding a new one after
# it, so back off by one character first.
$insPos = $lastPara.Range.End - 1
$insertionPoint = $d.Range($insPos, $insPos)

$newParagraphXml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">Im DM Do File </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
    <w:t>marital</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
    <w:t>status</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="de-DE"/>
    </w:rPr>
    <w:t xml:space="preserve"> korrigieren</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
